# Auto-generated edit script: update Leve profit calculations per scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1369.6888
$ws.Range("J112").Value = 1369.6888
$ws.Range("L112").Value = 4109.0664
$ws.Range("N112").Value = -6325.0664

$ws.Range("H132").Value = 21429.42
$ws.Range("I132").Value = 27440.48
$ws.Range("J132").Value = 5610.8423
$ws.Range("K132").Value = 82321.44
$ws.Range("L132").Value = 16832.5269
$ws.Range("M132").Value = -79791.44
$ws.Range("N132").Value = -21892.5269

$ws.Range("H137").Value = 868.71875
$ws.Range("I137").Value = 832.76
$ws.Range("J137").Value = 997.1429000000001
$ws.Range("K137").Value = 2498.28
$ws.Range("L137").Value = 2991.4287
$ws.Range("M137").Value = 51.72000000000025
$ws.Range("N137").Value = -8091.4287

$ws.Range("H138").Value = 3234.09
$ws.Range("I138").Value = 1454.36
$ws.Range("J138").Value = 5013.82
$ws.Range("K138").Value = 4363.08
$ws.Range("L138").Value = 15041.46
$ws.Range("M138").Value = 776.9200000000001
$ws.Range("N138").Value = -25321.46

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3037.25
$ws.Range("I32").Value = 3037.25
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3037.25
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -2750.25

$ws.Range("H74").Value = 1512.8918
$ws.Range("I74").Value = 1136.7273
$ws.Range("J74").Value = 2064.6
$ws.Range("K74").Value = 1136.7273
$ws.Range("L74").Value = 2064.6
$ws.Range("M74").Value = -262.7273
$ws.Range("N74").Value = -3812.6

$ws.Range("H77").Value = 1512.8918
$ws.Range("I77").Value = 1136.7273
$ws.Range("J77").Value = 2064.6
$ws.Range("K77").Value = 5683.636500000001
$ws.Range("L77").Value = 10323
$ws.Range("M77").Value = -1315.636500000001
$ws.Range("N77").Value = -19059

$ws.Range("H122").Value = 2273.48
$ws.Range("I122").Value = 1991.65
$ws.Range("K122").Value = 5974.950000000001
$ws.Range("M122").Value = -3524.950000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 579.3077
$ws.Range("I94").Value = 612.2
$ws.Range("J94").Value = 469.66666
$ws.Range("K94").Value = 612.2
$ws.Range("L94").Value = 469.66666
$ws.Range("M94").Value = -161.2
$ws.Range("N94").Value = -1371.66666

$ws.Range("H134").Value = 1486.6123
$ws.Range("I134").Value = 1122.8485
$ws.Range("K134").Value = 3368.5455
$ws.Range("M134").Value = -833.5455000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 46627.5
$ws.Range("J9").Value = 46627.5
$ws.Range("L9").Value = 46627.5
$ws.Range("N9").Value = -46963.5

$ws.Range("H31").Value = 3083.6187
$ws.Range("I31").Value = 1297.9851
$ws.Range("J31").Value = 7071.533
$ws.Range("K31").Value = 1297.9851
$ws.Range("L31").Value = 7071.533
$ws.Range("M31").Value = -1002.9851
$ws.Range("N31").Value = -7661.533

$ws.Range("H34").Value = 3083.6187
$ws.Range("I34").Value = 1297.9851
$ws.Range("J34").Value = 7071.533
$ws.Range("K34").Value = 1297.9851
$ws.Range("L34").Value = 7071.533
$ws.Range("M34").Value = -1095.9851
$ws.Range("N34").Value = -7475.533

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 15689.571
$ws.Range("I4").Value = 413.5
$ws.Range("J4").Value = 21800
$ws.Range("K4").Value = 1240.5
$ws.Range("L4").Value = 65400
$ws.Range("M4").Value = -1128.5
$ws.Range("N4").Value = -65624

$ws.Range("H56").Value = 5061.6665
$ws.Range("I56").Value = 5061.6665
$ws.Range("K56").Value = 5061.6665
$ws.Range("M56").Value = -4531.6665

$ws.Range("H107").Value = 426186.5
$ws.Range("I107").Value = 1399
$ws.Range("J107").Value = 648694.25
$ws.Range("K107").Value = 4197
$ws.Range("L107").Value = 1946082.75
$ws.Range("M107").Value = -2277
$ws.Range("N107").Value = -1949922.75

$ws.Range("H131").Value = 847.83
$ws.Range("J131").Value = 870.04254
$ws.Range("L131").Value = 2610.12762
$ws.Range("N131").Value = -12690.12762

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4149.75
$ws.Range("I5").Value = 1866.3334
$ws.Range("J5").Value = 11000
$ws.Range("K5").Value = 1866.3334
$ws.Range("L5").Value = 11000
$ws.Range("M5").Value = -1754.3334
$ws.Range("N5").Value = -11224

$ws.Range("H126").Value = 2033.9788
$ws.Range("I126").Value = 1849.5518
$ws.Range("J126").Value = 2331.111
$ws.Range("K126").Value = 5548.6554
$ws.Range("L126").Value = 6993.333
$ws.Range("M126").Value = -3078.6554
$ws.Range("N126").Value = -11933.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1046376.94
$ws.Range("I2").Value = 425375
$ws.Range("K2").Value = 425375
$ws.Range("M2").Value = -425263

$ws.Range("H7").Value = 62554.59
$ws.Range("I7").Value = 86719.086
$ws.Range("J7").Value = 4559.8
$ws.Range("K7").Value = 86719.086
$ws.Range("L7").Value = 4559.8
$ws.Range("M7").Value = -86607.086
$ws.Range("N7").Value = -4783.8

$ws.Range("H82").Value = 1816.6875
$ws.Range("I82").Value = 1137.1111
$ws.Range("J82").Value = 2690.4285
$ws.Range("K82").Value = 1137.1111
$ws.Range("L82").Value = 2690.4285
$ws.Range("M82").Value = -776.1111000000001
$ws.Range("N82").Value = -3412.4285

$ws.Range("H85").Value = 1816.6875
$ws.Range("I85").Value = 1137.1111
$ws.Range("J85").Value = 2690.4285
$ws.Range("K85").Value = 1137.1111
$ws.Range("L85").Value = 2690.4285
$ws.Range("M85").Value = 110.8888999999999
$ws.Range("N85").Value = -5186.4285

$ws.Range("H126").Value = 62554.59
$ws.Range("I126").Value = 86719.086
$ws.Range("J126").Value = 4559.8
$ws.Range("K126").Value = 260157.258
$ws.Range("L126").Value = 13679.4
$ws.Range("M126").Value = -257687.258
$ws.Range("N126").Value = -18619.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 22857.572
$ws.Range("J2").Value = 55001.5
$ws.Range("L2").Value = 55001.5
$ws.Range("N2").Value = -55225.5

$ws.Range("H132").Value = 1156.6
$ws.Range("I132").Value = 844.29785
$ws.Range("J132").Value = 2285.6924
$ws.Range("K132").Value = 2532.89355
$ws.Range("L132").Value = 6857.0772
$ws.Range("M132").Value = -2.893550000000232
$ws.Range("N132").Value = -11917.0772

$ws.Range("H136").Value = 5955593.5
$ws.Range("I136").Value = 9259771
$ws.Range("J136").Value = 8074
$ws.Range("K136").Value = 27779313
$ws.Range("L136").Value = 24222
$ws.Range("M136").Value = -27779313
$ws.Range("N136").Value = -29322
